# Completed entering 8PSK and 8APSK-L data
# Adds the 16APSK-L rows (13 new data rows) below the existing 16APSK block,
# and relocates the trailing "Short/Medium FECFRAME" summary rows further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the two trailing summary rows out of the way first (bottom-up) so
#    they are not overwritten while we populate rows 43-55.
$ws.Range("A47:C47").Cut($ws.Range("A64:C64"))
$ws.Range("A45:C45").Cut($ws.Range("A62:C62"))

# 2. New 16APSK-L rows - use the existing "16APSK" style block (row 35) as the
#    template for formatting (style s="2" on col A, s="1" on col B).
$ws.Range("A35:B35").Copy()
$ws.Range("A43:B55").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$labels = @(
    "16APSK-L 26/45",
    "16APSK-L 3/5",
    "16APSK-L 28/45",
    "16APSK-L 23/36",
    "16APSK-L 25/36",
    "16APSK-L 13/18",
    "16APSK-L 7/9",
    "16APSK-L 77/90",
    "16APSK-L 5/9",
    "16APSK-L 8/15",
    "16APSK-L 1/2",
    "16APSK-L 3/5",
    "16APSK-L 2/3"
)

$rates = @(
    0.57777777777777772,
    0.6,
    0.62222222222222223,
    0.63888888888888884,
    0.69444444444444442,
    0.72222222222222221,
    0.77777777777777779,
    0.85555555555555551,
    0.55555555555555558,
    0.53333333333333333,
    0.5,
    0.6,
    0.66666666666666663
)

for ($i = 0; $i -lt 13; $i++) {
    $r = 43 + $i
    $ws.Cells.Item($r, 1).Value = $labels[$i]
    $ws.Cells.Item($r, 2).Value = $rates[$i]
}

# 3. Sheet view: scroll back to the top and put the selection on G64.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("G64").Select()
